# Update the Prisoners_Data sheet: reorder the Category (A) / NumberOfPrisoners (B)
# values within several year-groups so that rows follow the F, C, M pattern.
# (Year and Country columns are untouched.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Year 2012 group (rows 5-7) ---
$ws.Range("A5").Value = "F"
$ws.Range("B5").Value = 11

$ws.Range("A7").Value = "M"
$ws.Range("B7").Value = 4656

# --- Year 2013 group (rows 8-10) ---
$ws.Range("A8").Value = "M"
$ws.Range("B8").Value = 5033

$ws.Range("A10").Value = "F"
$ws.Range("B10").Value = 16

# --- Year 2014 group (rows 11-13) ---
$ws.Range("A11").Value = "C"
$ws.Range("B11").Value = 156

$ws.Range("A12").Value = "M"
$ws.Range("B12").Value = 6200

$ws.Range("A13").Value = "F"
$ws.Range("B13").Value = 23

# --- Year 2017 group (rows 20-22) ---
$ws.Range("A20").Value = "F"
$ws.Range("B20").Value = 58

$ws.Range("A21").Value = "C"
$ws.Range("B21").Value = 350

$ws.Range("A22").Value = "M"
$ws.Range("B22").Value = 6171

# --- Year 2018 group (rows 23-25) ---
$ws.Range("A23").Value = "M"
$ws.Range("B23").Value = 5500

$ws.Range("A24").Value = "C"
$ws.Range("B24").Value = 230

$ws.Range("A25").Value = "F"
$ws.Range("B25").Value = 54

# --- Year 2020 group (rows 29-30) ---
$ws.Range("A29").Value = "C"
$ws.Range("B29").Value = 170

$ws.Range("A30").Value = "M"
$ws.Range("B30").Value = 4400
